# -----------------------------------------------------------------------
# Edit: rename header columns on "Weekly Quantity" / "Monthly Trend" sheets
#       and add a new "PO Forecast" sheet with forecast data (ds,
#       PO_Forecast, yhat_lower, yhat_upper).
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# 1) Rename the "Requested quantity" headers on the existing sheets.
$wsWeekly  = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

$wsWeekly.Range("B1").Value  = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# 2) Add the new "PO Forecast" worksheet as the last (3rd) tab.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row.
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"
$wsForecast.Range("A1:D1").Font.Bold = $true

$poForecastData = @(
    ,@(45025.99999999999, 71, -28.21677656631963, 162.3380479894555)
    ,@(45032.99999999999, 71, -20.23696643990653, 167.2507691192567)
    ,@(45039.99999999999, 71, -17.65701230870259, 169.0397230316897)
    ,@(45046.99999999999, 71, -20.52268899855125, 168.6524718851996)
    ,@(45060.99999999999, 71, -24.56225615910731, 165.5650560364483)
    ,@(45067.99999999999, 72, -26.33268383353673, 164.8034595619595)
    ,@(45088.99999999999, 72, -21.91007678836576, 165.1702328434654)
    ,@(45102.99999999999, 72, -27.6447699834739, 158.1506914757409)
    ,@(45109.99999999999, 72, -17.43974685344936, 163.8526918462009)
    ,@(45137.99999999999, 73, -20.12674799706491, 171.9751416287586)
    ,@(45144.99999999999, 73, -20.89572681579227, 168.3745591953663)
    ,@(45151.99999999999, 73, -13.83989929306486, 172.2474947698102)
    ,@(45158.99999999999, 73, -25.18487679769967, 167.6364626559885)
    ,@(45179.99999999999, 74, -15.84022129092052, 170.7378518674577)
    ,@(45193.99999999999, 74, -16.14166329857698, 165.8402399052939)
    ,@(45200.99999999999, 74, -22.73782952755248, 164.76672584292)
    ,@(45207.99999999999, 74, -20.61905200321639, 165.1278724870119)
    ,@(45214.99999999999, 74, -25.66670362764065, 169.617248351972)
    ,@(45221.99999999999, 74, -19.44818773672649, 171.7660428424636)
    ,@(45235.99999999999, 75, -30.51098614521884, 167.2791627419558)
    ,@(45242.99999999999, 75, -18.20192904159948, 171.9218348819983)
    ,@(45249.99999999999, 75, -20.90774713914364, 169.4920980919215)
    ,@(45256.99999999999, 75, -18.24229267485949, 174.5556780555137)
    ,@(45263.99999999999, 75, -20.06262939179571, 175.3438547878284)
    ,@(45305.99999999999, 76, -22.21200404089453, 166.2545909491403)
    ,@(45312.99999999999, 76, -19.44838767845842, 177.7713985188225)
    ,@(45333.99999999999, 76, -13.90005637103124, 177.2900872365813)
    ,@(45340.99999999999, 76, -12.89329448682815, 178.818643085625)
    ,@(45347.99999999999, 77, -15.23602338176387, 170.6990546768848)
    ,@(45354.99999999999, 77, -15.41129391876943, 172.5515145491383)
    ,@(45361.99999999999, 77, -19.91829413557062, 166.0902409854458)
    ,@(45368.99999999999, 77, -19.53532880670817, 170.2588158024665)
    ,@(45375.99999999999, 77, -11.24637062832593, 172.9983723912783)
    ,@(45382.99999999999, 77, -19.26448950181715, 174.682794921398)
    ,@(45389.99999999999, 77, -23.03594471694328, 175.8704527150077)
    ,@(45396.99999999999, 77, -21.0073518010217, 168.8330297724889)
    ,@(45410.99999999999, 78, -17.91151757654313, 166.3112498922571)
    ,@(45424.99999999999, 78, -15.97809053185206, 172.1447989404428)
    ,@(45431.99999999999, 78, -19.36042060796551, 164.4156176725316)
    ,@(45438.99999999999, 78, -23.25298230668656, 168.2139319383576)
    ,@(45445.99999999999, 78, -15.70703188189221, 170.8734981069626)
    ,@(45459.99999999999, 79, -21.77249078720494, 169.8858415502849)
    ,@(45466.99999999999, 79, -20.44095025836367, 170.0035658773764)
    ,@(45487.99999999999, 79, -12.0758054794917, 172.8127403795625)
    ,@(45501.99999999999, 79, -12.03407418557256, 169.9662384416014)
    ,@(45508.99999999999, 79, -12.63844218422951, 168.8848454291002)
    ,@(45515.99999999999, 80, -17.59024481947676, 170.4264850595441)
    ,@(45536.99999999999, 80, -16.79673398350523, 174.6669229697923)
    ,@(45550.99999999999, 80, -11.36877780055335, 180.7462502343483)
    ,@(45557.99999999999, 80, -8.647722560036133, 181.7845080466439)
    ,@(45578.99999999999, 81, -7.08780310255937, 183.0035622268193)
    ,@(45592.99999999999, 81, -14.74292080712049, 176.0512794096134)
    ,@(45599.99999999999, 81, -12.25300056397103, 165.8645057427096)
    ,@(45606.99999999999, 81, -12.13414849823725, 177.0562557809844)
    ,@(45613.99999999999, 81, -7.09780826831974, 186.3239439613612)
    ,@(45620.99999999999, 81, -9.399311447932773, 180.6110862070998)
    ,@(45627.99999999999, 82, -12.10204201869205, 177.8053569643739)
    ,@(45634.99999999999, 82, -13.5502865914628, 182.5578556037352)
    ,@(45641.99999999999, 82, -12.53776903144552, 175.5038197138916)
    ,@(45648.99999999999, 82, -7.691837888715024, 172.8704516323194)
    ,@(45655.99999999999, 82, -4.420359950316651, 179.9877705915623)
)

# 3) Write the forecast rows starting at row 2, mirroring the date-format
#    style used for column A ("ds") on the other two sheets.
for ($i = 0; $i -lt $poForecastData.Count; $i++) {
    $r = $i + 2
    $row = $poForecastData[$i]

    $dsCell = $wsForecast.Cells.Item($r, 1)
    $dsCell.Value = $row[0]
    $dsCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
}

Write-Host "PO Forecast sheet populated with $($poForecastData.Count) data rows."
